$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing OT values in E2:E4 (plain alnum text, no risk of numeric coercion) ---
$ws.Range("E2").Value = 'ICD30466266'
$ws.Range("E3").Value = 'ICD30466317'
$ws.Range("E4").Value = 'ICD30465943'

# --- Append new rows 27-30 ---
# Row 27
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = '4966'
$ws.Range("A27").ClearFormats()
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '8/22/2025'
$ws.Range("B27").ClearFormats()
$ws.Range("C27").Value = 'MUÑECAS 1035'
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = 'ICD30463961'
$ws.Range("F27").Value = 'Optical Power'
$ws.Range("G27").Value = 'Pendiente'
$ws.Range("H27").Value = 'Caja de empalme colgando y cable cortado'
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = '{"direccionesNormalizadas": [{"altura": 1035, "cod_calle": 13141, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.448024", "y": "-34.598462"}, "direccion": "MUÑECAS 1035, CABA", "nombre_calle": "MUÑECAS", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K27").Value = -58.448024
$ws.Range("L27").Value = -34.598462
$ws.Range("M27").Value = 'Paternal'
$ws.Range("N27").Value = 'Capital Norte'

# Row 28
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = '5117'
$ws.Range("A28").ClearFormats()
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = '8/22/2025'
$ws.Range("B28").ClearFormats()
$ws.Range("C28").Value = 'ARENAL, CONCEPCION 3470'
$ws.Range("D28").Value = 15
$ws.Range("E28").Value = 'ICD30465386'
$ws.Range("F28").Value = 'Optical Power'
$ws.Range("G28").Value = 'Pendiente'
$ws.Range("H28").Value = 'Cable en panza y cortados'
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = '{"direccionesNormalizadas": [{"altura": 3470, "cod_calle": 1103, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.446142", "y": "-34.583381"}, "direccion": "ARENAL, CONCEPCION 3470, CABA", "nombre_calle": "ARENAL, CONCEPCION", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K28").Value = -58.446142
$ws.Range("L28").Value = -34.583381
$ws.Range("M28").Value = 'Colegiales'
$ws.Range("N28").Value = 'Capital Norte'

# Row 29
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = '5508'
$ws.Range("A29").ClearFormats()
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = '8/22/2025'
$ws.Range("B29").ClearFormats()
$ws.Range("C29").Value = 'BRUNO, GIORDANO 829'
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 'ICD30463299'
$ws.Range("F29").Value = 'Optical Power'
$ws.Range("G29").Value = 'Pendiente'
$ws.Range("H29").Value = 'Tendido aereo en panza'
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = '{"direccionesNormalizadas": [{"altura": 829, "cod_calle": 2118, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.441705", "y": "-34.617573"}, "direccion": "BRUNO, GIORDANO 829, CABA", "nombre_calle": "BRUNO, GIORDANO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K29").Value = -58.441705
$ws.Range("L29").Value = -34.617573
$ws.Range("M29").Value = 'Almagro'
$ws.Range("N29").Value = 'Capital Sur'

# Row 30
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = '5571'
$ws.Range("A30").ClearFormats()
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = '8/22/2025'
$ws.Range("B30").ClearFormats()
$ws.Range("C30").Value = 'ESTADO DE ISRAEL AV. 4624'
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 'ICD30463349'
$ws.Range("F30").Value = 'Optical Power'
$ws.Range("G30").Value = 'Pendiente'
$ws.Range("H30").Value = 'Tendido a baja altura'
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = '{"direccionesNormalizadas": [{"altura": 4624, "cod_calle": 5086, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.430250", "y": "-34.600000"}, "direccion": "ESTADO DE ISRAEL AV. 4624, CABA", "nombre_calle": "ESTADO DE ISRAEL AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K30").Value = -58.43025
$ws.Range("L30").Value = -34.6
$ws.Range("M30").Value = 'Almagro'
$ws.Range("N30").Value = 'Capital Sur'
